$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) values for the cryptos list

$ws.Cells.Item(2, 4).Value = '41.523.92'
$ws.Cells.Item(2, 5).Value = '  +0.07%  '
$ws.Cells.Item(3, 4).Value = '2.463.55'
$ws.Cells.Item(3, 5).Value = '  -0.63%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.54%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '313.45'
$ws.Cells.Item(5, 5).Value = '  -0.04%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '91.23'
$ws.Cells.Item(6, 5).Value = '  -1.86%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.549'
$ws.Cells.Item(7, 5).Value = '  +0.51%  '
$ws.Cells.Item(8, 5).Value = '  -0.55%  '
$ws.Cells.Item(9, 5).Value = '  +3.78%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '32.43'
$ws.Cells.Item(10, 5).Value = '  -2.29%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0791'
$ws.Cells.Item(11, 5).Value = '  +1.50%  '
$ws.Cells.Item(12, 5).Value = '  +0.49%  '
$ws.Cells.Item(13, 5).Value = '  -0.53%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.83'
$ws.Cells.Item(14, 5).Value = '  -0.57%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '15.76'
$ws.Cells.Item(15, 5).Value = '  +2.29%  '
$ws.Cells.Item(16, 4).Value = '2.454.29'
$ws.Cells.Item(16, 5).Value = '  -2.96%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.774'
$ws.Cells.Item(17, 5).Value = '  -1.44%  '
$ws.Cells.Item(18, 4).Value = '41.491.59'
$ws.Cells.Item(18, 5).Value = '  +0.33%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0940'
$ws.Cells.Item(20, 5).Value = '  +1.75%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '70.82'
$ws.Cells.Item(21, 5).Value = '  +1.29%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '11.06'
$ws.Cells.Item(22, 5).Value = '  -0.52%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '237.20'
$ws.Cells.Item(23, 5).Value = '  +1.09%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.72'
$ws.Cells.Item(24, 5).Value = '  -0.96%  '
$ws.Cells.Item(25, 5).Value = '  +0.01%  '
$ws.Cells.Item(26, 5).Value = '  +1.09%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '24.47'
$ws.Cells.Item(27, 5).Value = '  +1.84%  '
$ws.Cells.Item(28, 5).Value = '  -0.35%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.67'
$ws.Cells.Item(29, 5).Value = '  -0.95%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '35.24'
$ws.Cells.Item(30, 5).Value = '  -3.09%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '155.77'
$ws.Cells.Item(31, 5).Value = '  +1.84%  '
$ws.Cells.Item(32, 5).Value = '  -0.65%  '
$ws.Cells.Item(33, 5).Value = '  +0.74%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0757'
$ws.Cells.Item(34, 5).Value = '  +1.08%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '17.08'
$ws.Cells.Item(35, 5).Value = '  -3.31%  '
$ws.Cells.Item(36, 5).Value = '  -7.05%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.87'
$ws.Cells.Item(37, 5).Value = '  -5.35%  '
$ws.Cells.Item(38, 5).Value = '  +1.02%  '
$ws.Cells.Item(39, 5).Value = '  +2.50%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.78'
$ws.Cells.Item(40, 5).Value = '  -4.13%  '
$ws.Cells.Item(41, 5).Value = '  -1.89%  '
$ws.Cells.Item(42, 5).Value = '  -0.80%  '
$ws.Cells.Item(43, 4).Value = '1.941.26'
$ws.Cells.Item(43, 5).Value = '  -1.52%  '
$ws.Cells.Item(44, 5).Value = '  -0.03%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '18.74'
$ws.Cells.Item(45, 5).Value = '  -4.73%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.89'
$ws.Cells.Item(46, 5).Value = '  -2.61%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '9.01'
$ws.Cells.Item(47, 5).Value = '  +2.73%  '
$ws.Cells.Item(48, 4).Value = '2.706.42'
$ws.Cells.Item(48, 5).Value = '  -0.68%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '96.55'
$ws.Cells.Item(49, 5).Value = '  +0.42%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '66.83'
$ws.Cells.Item(50, 5).Value = '  -2.11%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '52.02'
$ws.Cells.Item(51, 5).Value = '  +3.47%  '
